$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Row 12: ProfileSummaryTest / Add summary validation / Y / (blank) ---
$ws.Range("B12").Value = "Add summary validation"
$ws.Range("A12").Value = "ProfileSummaryTest"
$ws.Range("C12").Value = "Y"
$ws.Range("D12").Value = ""

# --- Row 13: ProfileFollowerTest / To verify count of users following me / Y / (blank) ---
$ws.Range("A13").Value = "ProfileFollowerTest"
$ws.Range("B13").Value = "To verify count of users following me"
$ws.Range("C13").Value = "Y"
$ws.Range("D13").Value = ""

# --- Copy formatting (thin border, no fill) from existing rows onto the new rows ---
$ws.Range("A5").Copy()
$ws.Range("A12").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("B5").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B13").PasteSpecial(-4122)

$ws.Range("D5").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("D13").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C13").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Match the author's final selection state ---
$ws.Range("J11").Select()
